$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header label in B1 (was "nom_w_pc_obs")
$ws.Range("B1").Value = "govdebt_rcpc_obs"

# Update B2:B103 with revised series values
$values = @(
    0.1261408116688424,
    0.09353756063987154,
    0.09384950029112429,
    0.07578909220622135,
    0.05384392728215217,
    0.05490798724647763,
    0.05213786959980706,
    0.0555075466825802,
    0.03160103075550027,
    0.05133398833675931,
    0.03939551562873391,
    0.0364179876403933,
    0.03007042301872578,
    0.03313772103084092,
    0.03302719713448945,
    0.01884574300483135,
    0.01684175144963231,
    0.007662325536864543,
    0.006539888365702538,
    0.007310038471498047,
    0.003609078723670799,
    0.0001134386481892666,
    -0.005651770378015321,
    -0.00119886424152358,
    -0.008269216071823272,
    -0.005154515268939931,
    -0.001238423165013415,
    -0.001409740693218843,
    -0.009179236696036021,
    -0.003812588617461678,
    -0.002726727613686065,
    0.002957800509127438,
    -0.007221190532190704,
    0.003434289841226644,
    0.005695085355889383,
    0.006015094549419196,
    0.009855983884072515,
    0.00922330859072109,
    0.00995093723119269,
    0.005487405253851073,
    0.00542943473688564,
    0.001406276967416041,
    0.001601895754599035,
    -0.002750166649822951,
    -0.005567676091189667,
    -0.007767653769757167,
    -0.007188200741892949,
    -0.007727092537029946,
    -0.007204523533620978,
    -0.008174065050497205,
    -0.01069349520463034,
    -0.01186743252802933,
    -0.01052187789044488,
    -0.01307135461899232,
    -0.0159552500397838,
    -0.01936573506249615,
    -0.02225760244998063,
    -0.02262002628295998,
    -0.02192668707309529,
    -0.02295937330100038,
    -0.02421747329354793,
    -0.02680865738743186,
    -0.02997475873340192,
    -0.02916003705531565,
    -0.03372286220034435,
    -0.03333563999300391,
    -0.03319222732820069,
    -0.03368925618571903,
    -0.05652797037467869,
    -0.04117613723731437,
    -0.04262647510796067,
    -0.04158153678582131,
    -0.04181908127961918,
    -0.038596037427966,
    -0.01598660192597055,
    -0.02335102372976414,
    -0.007228579935758851,
    -0.004529302573401985,
    -0.004294746697896513,
    -0.00277811403998416,
    -0.008614446755511686,
    0.005705856502581284,
    0.007383173408631948,
    0.000661995384704149,
    0.001856653328962253,
    -0.003242639865522846,
    -0.00363266611469313,
    -0.005093233057539372,
    -0.01088737402535443,
    -0.008998977569558121,
    -0.009610680291573061,
    -0.01303141121175187,
    -0.01464895280416232,
    -0.01540808523565693,
    -0.01548267581010285,
    -0.01816524297492648,
    -0.02107375400066789,
    -0.01318356937406171,
    -0.008231480558097949,
    -0.009630314692508392,
    -0.002522855657092073,
    0.01124975073282699
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}

# Row 104 no longer has a value in column B; clear it
$ws.Range("B104").ClearContents()

Write-Host "Applied data update for govdebt_rcpc_obs series"